$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: replicate the formatting of column K onto the corresponding column L
# cell (font, color mode, borders, alignment, number format), matching the
# style groups used by rows 4-18 of this worksheet.
function Format-LikeColumnK {
    param(
        [string]$Addr,
        [bool]$Bold,
        [string]$ColorMode,
        [bool]$TopBorder,
        [bool]$BottomBorder,
        [string]$Fmt
    )
    $r = $ws.Range($Addr)
    $r.Font.Name = "Times New Roman"
    $r.Font.Size = 9
    $r.Font.Bold = $Bold
    if ($ColorMode -eq "theme") {
        $r.Font.ThemeColor = 1
    } else {
        $r.Font.Color = 0
    }
    if ($TopBorder) {
        $r.Borders.Item(8).LineStyle = 1
        $r.Borders.Item(8).Weight = -4138
    }
    if ($BottomBorder) {
        $r.Borders.Item(9).LineStyle = 1
        $r.Borders.Item(9).Weight = -4138
    }
    $r.HorizontalAlignment = -4152
    $r.VerticalAlignment = -4108
    $r.WrapText = $true
    if ($Fmt -ne "") {
        $r.NumberFormat = $Fmt
    }
}

# Row 4: new year header column (2022)
$ws.Range("L4").Value = 2022
Format-LikeColumnK "L4" $true "theme" $true $true ""

# Row 5: total
$ws.Range("L5").Value = 8800.6
Format-LikeColumnK "L5" $true "rgb" $false $false "0.0"

# Row 6: section header (empty data cell)
Format-LikeColumnK "L6" $false "rgb" $false $false "0.0"

# Row 7: formula L5-L8
$ws.Range("L7").Formula = "=L5-L8"
Format-LikeColumnK "L7" $false "rgb" $false $false "0.0"

# Row 8
$ws.Range("L8").Value = 258.39999999999998
Format-LikeColumnK "L8" $false "rgb" $false $false "0.0"

# Row 9: section header (empty data cell)
Format-LikeColumnK "L9" $true "theme" $false $false "0.0"

# Row 10
$ws.Range("L10").Value = 683.8
Format-LikeColumnK "L10" $false "theme" $false $false "0.0"

# Row 11
$ws.Range("L11").Value = 1101.8
Format-LikeColumnK "L11" $false "theme" $false $false "0.0"

# Row 12
$ws.Range("L12").Value = 714.9
Format-LikeColumnK "L12" $false "theme" $false $false "0.0"

# Row 13
$ws.Range("L13").Value = 757.9
Format-LikeColumnK "L13" $false "theme" $false $false "0.0"

# Row 14
$ws.Range("L14").Value = 1383.3
Format-LikeColumnK "L14" $false "theme" $false $false "0.0"

# Row 15
$ws.Range("L15").Value = 1023.7
Format-LikeColumnK "L15" $false "theme" $false $false "0.0"

# Row 16
$ws.Range("L16").Value = 2929.3
Format-LikeColumnK "L16" $false "theme" $false $false "0.0"

# Row 17
$ws.Range("L17").Value = 148.9
Format-LikeColumnK "L17" $false "theme" $false $false "0.0"

# Row 18: totals row (bottom medium border)
$ws.Range("L18").Value = 57
Format-LikeColumnK "L18" $false "theme" $false $true "0.0"

# Update the selected cell to match the authored workbook state
$null = $ws.Range("M4").Select()
